# Applies the crypto price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of new Price values are plain numeric-looking strings (e.g. "298.53").
# The sheet stores the whole Price column as text, so - just like a user would have to
# do by hand in the Excel UI - we switch those specific cells to the Text number format
# before writing them; otherwise Excel would silently convert them to numbers.
$textForcedCells = @(
    "D5", "D6", "D9", "D11", "D12", "D13", "D16", "D17",
    "D19", "D21", "D22", "D23", "D25", "D26", "D28", "D29",
    "D30", "D31", "D32", "D34", "D35", "D36", "D39", "D40",
    "D41", "D43", "D46", "D51"
)
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).NumberFormat = "@"
}


# Row 2
$ws.Range('D2').Value = '43.134.20'
$ws.Range('E2').Value = '  -6.32%  '

# Row 3
$ws.Range('D3').Value = '2.550.29'
$ws.Range('E3').Value = '  -2.46%  '

# Row 4
$ws.Range('E4').Value = '  -0.19%  '

# Row 5
$ws.Range('D5').Value = '298.53'
$ws.Range('E5').Value = '  -3.64%  '

# Row 6
$ws.Range('D6').Value = '94.38'
$ws.Range('E6').Value = '  -4.79%  '

# Row 7
$ws.Range('E7').Value = '  -3.41%  '

# Row 8
$ws.Range('E8').Value = '  +0.01%  '

# Row 9
$ws.Range('D9').Value = '0.552'
$ws.Range('E9').Value = '  -5.00%  '

# Row 10
$ws.Range('E10').Value = '  -7.30%  '

# Row 11
$ws.Range('D11').Value = '0.0811'
$ws.Range('E11').Value = '  -4.04%  '

# Row 12
$ws.Range('D12').Value = '7.76'
$ws.Range('E12').Value = '  -3.94%  '

# Row 13
$ws.Range('D13').Value = '0.109'
$ws.Range('E13').Value = '  +1.44%  '

# Row 14
$ws.Range('D14').Value = '2.941.98'
$ws.Range('E14').Value = '  -2.50%  '

# Row 15
$ws.Range('D15').Value = '2.536.20'
$ws.Range('E15').Value = '  -3.00%  '

# Row 16
$ws.Range('D16').Value = '0.874'
$ws.Range('E16').Value = '  -4.44%  '

# Row 17
$ws.Range('D17').Value = '14.21'
$ws.Range('E17').Value = '  -4.23%  '

# Row 18
$ws.Range('D18').Value = '43.146.28'
$ws.Range('E18').Value = '  -6.81%  '

# Row 19
$ws.Range('D19').Value = '13.03'
$ws.Range('E19').Value = '  +2.31%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0984'
$ws.Range('E20').Value = '  -3.15%  '

# Row 21
$ws.Range('D21').Value = '6.66'
$ws.Range('E21').Value = '  -1.49%  '

# Row 22
$ws.Range('D22').Value = '72.13'
$ws.Range('E22').Value = '  -1.41%  '

# Row 23
$ws.Range('D23').Value = '260.72'
$ws.Range('E23').Value = '  -10.46%  '

# Row 24
$ws.Range('E24').Value = '  -4.14%  '

# Row 25
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').Value = '2.15'
$ws.Range('E25').Value = '  -4.97%  '

# Row 26
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = '29.41'
$ws.Range('E26').Value = '  -0.72%  '

# Row 27
$ws.Range('E27').Value = '  +0.28%  '

# Row 28
$ws.Range('D28').Value = '10.07'
$ws.Range('E28').Value = '  -6.91%  '

# Row 29
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').Value = '36.99'
$ws.Range('E29').Value = '  -4.99%  '

# Row 30
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '2.13'
$ws.Range('E30').Value = '  -3.94%  '

# Row 31
$ws.Range('D31').Value = '6.01'
$ws.Range('E31').Value = '  -4.21%  '

# Row 32
$ws.Range('D32').Value = '155.25'
$ws.Range('E32').Value = '  -2.21%  '

# Row 33
$ws.Range('E33').Value = '  -3.72%  '

# Row 34
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').Value = '2.74'
$ws.Range('E34').Value = '  -2.17%  '

# Row 35
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = '3.39'
$ws.Range('E35').Value = '  -6.04%  '

# Row 36
$ws.Range('D36').Value = '0.0801'
$ws.Range('E36').Value = '  -4.93%  '

# Row 37
$ws.Range('E37').Value = '  -5.39%  '

# Row 38
$ws.Range('E38').Value = '  -2.92%  '

# Row 39
$ws.Range('D39').Value = '16.58'
$ws.Range('E39').Value = '  +5.61%  '

# Row 40
$ws.Range('D40').Value = '23.29'
$ws.Range('E40').Value = '  +7.13%  '

# Row 41
$ws.Range('D41').Value = '3.49'
$ws.Range('E41').Value = '  -2.18%  '

# Row 42
$ws.Range('E42').Value = '  -5.14%  '

# Row 43
$ws.Range('D43').Value = '3.90'
$ws.Range('E43').Value = '  -3.39%  '

# Row 44
$ws.Range('D44').Value = '2.069.00'
$ws.Range('E44').Value = '  -2.82%  '

# Row 45
$ws.Range('E45').Value = '  -0.18%  '

# Row 46
$ws.Range('D46').Value = '86.02'
$ws.Range('E46').Value = '  -11.09%  '

# Row 47
$ws.Range('E47').Value = '  +2.74%  '

# Row 48
$ws.Range('D48').Value = '2.797.30'
$ws.Range('E48').Value = '  -2.66%  '

# Row 49
$ws.Range('E49').Value = '  -7.86%  '

# Row 50
$ws.Range('E50').Value = '  -1.93%  '

# Row 51
$ws.Range('D51').Value = '104.36'
$ws.Range('E51').Value = '  -5.04%  '
